$d = $word.ActiveDocument

$pairs = @(
    @("548×2=", "700×7="),
    @("815×6=", "583×3="),
    @("983×9=", "528×2="),
    @("930×6=", "997×5="),
    @("307×3=", "122×8="),
    @("203×8=", "708×3="),
    @("261×8=", "612×2="),
    @("126×7=", "991×9="),
    @("982×6=", "489×9="),
    @("686×2=", "987×8="),
    @("279×7=", "658×2="),
    @("941×6=", "602×6="),
    @("376×5=", "765×6="),
    @("718×9=", "260×9="),
    @("932×8=", "417×8="),
    @("933×8=", "432×3="),
    @("403×7=", "286×3="),
    @("578×9=", "114×6="),
    @("788×4=", "340×7="),
    @("248×5=", "381×5="),
    @("314×6=", "130×2="),
    @("344×9=", "192×8="),
    @("200×4=", "594×3="),
    @("376×3=", "452×9="),
    @("312×3=", "411×2="),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
